$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.049372413810607
$ws.Range("D2").Value = 1.051811990725495
$ws.Range("E2").Value = 1.046350559679824
$ws.Range("F2").Value = 1.052953261296608
$ws.Range("I2").Value = 1.048203490450746
$ws.Range("J2").Value = 1.054411209305147
$ws.Range("K2").Value = 1.054562456638632
$ws.Range("L2").Value = 1.049116240933994
$ws.Range("M2").Value = 1.055700569281105
$ws.Range("N2").Value = 1.055908594165016
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.050728646042662
$ws.Range("D3").Value = 1.052915471161338
$ws.Range("E3").Value = 1.047544295860482
$ws.Range("F3").Value = 1.055430357076857
$ws.Range("I3").Value = 1.048774040477393
$ws.Range("J3").Value = 1.055414704330077
$ws.Range("K3").Value = 1.055477910815659
$ws.Range("L3").Value = 1.050120609141667
$ws.Range("M3").Value = 1.057986353182528
$ws.Range("N3").Value = 1.056913514267985
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.051603642098911
$ws.Range("D4").Value = 1.053627166168802
$ws.Range("E4").Value = 1.048313904880858
$ws.Range("F4").Value = 1.05702596693558
$ws.Range("I4").Value = 1.049140225469716
$ws.Range("J4").Value = 1.056061003110874
$ws.Range("K4").Value = 1.05606735593772
$ws.Range("L4").Value = 1.050767122708032
$ws.Range("M4").Value = 1.059457896722576
$ws.Range("N4").Value = 1.057560730867182
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.051970883785694
$ws.Range("D5").Value = 1.053925813017544
$ws.Range("E5").Value = 1.04863678300665
$ws.Range("F5").Value = 1.057695073177201
$ws.Range("I5").Value = 1.04929345812808
$ws.Range("J5").Value = 1.05633199031778
$ws.Range("K5").Value = 1.056314468729562
$ws.Range("L5").Value = 1.051038116544109
$ws.Range("M5").Value = 1.060074777025021
$ws.Range("N5").Value = 1.057832102907004
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.052032510002125
$ws.Range("D6").Value = 1.053975925137947
$ws.Range("E6").Value = 1.048690956903159
$ws.Range("F6").Value = 1.057807321341825
$ws.Range("I6").Value = 1.049319145028308
$ws.Range("J6").Value = 1.056377448473504
$ws.Range("K6").Value = 1.056355919824229
$ws.Range("L6").Value = 1.051083570911877
$ws.Range("M6").Value = 1.060178252132865
$ws.Range("N6").Value = 1.057877625618524
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.051608551569922
$ws.Range("D7").Value = 1.053631158851349
$ws.Range("E7").Value = 1.048318221796539
$ws.Range("F7").Value = 1.05703491413985
$ws.Range("I7").Value = 1.049142275758921
$ws.Range("J7").Value = 1.05606462686005
$ws.Range("K7").Value = 1.056070660571051
$ws.Range("L7").Value = 1.050770746874287
$ws.Range("M7").Value = 1.059466146353863
$ws.Range("N7").Value = 1.057564359762497
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.04983129861273
$ws.Range("D8").Value = 1.052185404656841
$ws.Range("E8").Value = 1.046754577461223
$ws.Range("F8").Value = 1.053791932718614
$ws.Range("I8").Value = 1.048396935386319
$ws.Range("J8").Value = 1.054750977588506
$ws.Range("K8").Value = 1.054872447011157
$ws.Range("L8").Value = 1.049456376934805
$ws.Range("M8").Value = 1.056474641217383
$ws.Range("N8").Value = 1.056248844958309
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.046679352004719
$ws.Range("D9").Value = 1.049619578628892
$ws.Range("E9").Value = 1.04397723881623
$ws.Range("F9").Value = 1.048019947905993
$ws.Range("I9").Value = 1.047060288529742
$ws.Range("J9").Value = 1.052412572089623
$ws.Range("K9").Value = 1.052738358430352
$ws.Range("L9").Value = 1.04711401153923
$ws.Range("M9").Value = 1.051143807261036
$ws.Range("N9").Value = 1.053907118655399
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.044563818200775
$ws.Range("D10").Value = 1.047896249136287
$ws.Range("E10").Value = 1.042110305043727
$ws.Range("F10").Value = 1.044130498630351
$ws.Range("I10").Value = 1.046153136385
$ws.Range("J10").Value = 1.050837232847835
$ws.Range("K10").Value = 1.051299875407612
$ws.Range("L10").Value = 1.045534213337062
$ws.Range("M10").Value = 1.047547297206493
$ws.Range("N10").Value = 1.052329542251195
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.043644244582951
$ws.Range("D11").Value = 1.047146878402603
$ws.Range("E11").Value = 1.041298123647158
$ws.Range("F11").Value = 1.042435849922837
$ws.Range("I11").Value = 1.045756428069454
$ws.Range("J11").Value = 1.050151080627727
$ws.Range("K11").Value = 1.050673144625614
$ws.Range("L11").Value = 1.044845696664044
$ws.Range("M11").Value = 1.045979257546225
$ws.Range("N11").Value = 1.051642415616231
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.043302130261325
$ws.Range("D12").Value = 1.046868043927556
$ws.Range("E12").Value = 1.040995862892096
$ws.Range("F12").Value = 1.041804750029901
$ws.Range("I12").Value = 1.045608478001304
$ws.Range("J12").Value = 1.049895598485941
$ws.Range("K12").Value = 1.050439759269553
$ws.Range("L12").Value = 1.044589270911723
$ws.Range("M12").Value = 1.04539515383777
$ws.Range("N12").Value = 1.0513865706605
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.043375539758373
$ws.Range("D13").Value = 1.046927876976969
$ws.Range("E13").Value = 1.041060725301918
$ws.Range("F13").Value = 1.041940197851756
$ws.Range("I13").Value = 1.04564024084815
$ws.Range("J13").Value = 1.049950428260734
$ws.Range("K13").Value = 1.050489848065221
$ws.Range("L13").Value = 1.044644306054287
$ws.Range("M13").Value = 1.045520522214696
$ws.Range("N13").Value = 1.051441478299861
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.04361597647354
$ws.Range("D14").Value = 1.047123839808896
$ws.Range("E14").Value = 1.04127315058644
$ws.Range("F14").Value = 1.042383716625268
$ws.Range("I14").Value = 1.04574421065504
$ws.Range("J14").Value = 1.050129974995516
$ws.Range("K14").Value = 1.050653865006119
$ws.Range("L14").Value = 1.044824514366838
$ws.Range("M14").Value = 1.045931009602248
$ws.Range("N14").Value = 1.051621280011602
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.043764045030816
$ws.Range("D15").Value = 1.047244514419716
$ws.Range("E15").Value = 1.041403955556557
$ws.Range("F15").Value = 1.042656765251542
$ws.Range("I15").Value = 1.045808190802401
$ws.Range("J15").Value = 1.050240517893567
$ws.Range("K15").Value = 1.050754842822064
$ws.Range("L15").Value = 1.044935456220029
$ws.Range("M15").Value = 1.046183702355129
$ws.Range("N15").Value = 1.051731979893247
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.044624771732231
$ws.Range("D16").Value = 1.047945914990063
$ws.Range("E16").Value = 1.042164126087391
$ws.Range("F16").Value = 1.044242741111859
$ws.Range("I16").Value = 1.046179381630363
$ws.Range("J16").Value = 1.050882684905357
$ws.Range("K16").Value = 1.051341387329276
$ws.Range("L16").Value = 1.045579813161783
$ws.Range("M16").Value = 1.047651132382846
$ws.Range("N16").Value = 1.05237505885585
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.045163727804489
$ws.Range("D17").Value = 1.048385032335022
$ws.Range("E17").Value = 1.042639939086268
$ws.Range("F17").Value = 1.045234734024094
$ws.Range("I17").Value = 1.046411168493272
$ws.Range("J17").Value = 1.051284415548064
$ws.Range("K17").Value = 1.05170827126062
$ws.Range("L17").Value = 1.045982801241911
$ws.Range("M17").Value = 1.048568703894442
$ws.Range("N17").Value = 1.052777360002151
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.045477751708429
$ws.Range("D18").Value = 1.048640858093849
$ws.Range("E18").Value = 1.042917108177668
$ws.Range("F18").Value = 1.045812338242206
$ws.Range("I18").Value = 1.046545989589657
$ws.Range("J18").Value = 1.051518351101882
$ws.Range("K18").Value = 1.051921896879594
$ws.Range("L18").Value = 1.04621742840508
$ws.Range("M18").Value = 1.049102876411102
$ws.Range("N18").Value = 1.053011627771288
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.045584768451672
$ws.Range("D19").Value = 1.048728036822066
$ws.Range("E19").Value = 1.043011554181364
$ws.Range("F19").Value = 1.046009116970202
$ws.Range("I19").Value = 1.046591896524918
$ws.Range("J19").Value = 1.05159805173713
$ws.Range("K19").Value = 1.051994674970147
$ws.Range("L19").Value = 1.046297357771536
$ws.Range("M19").Value = 1.049284842057156
$ws.Range("N19").Value = 1.053091441590579
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.045105938184772
$ws.Range("D20").Value = 1.048337950714638
$ws.Range("E20").Value = 1.042588926647804
$ws.Range("F20").Value = 1.045128407282635
$ws.Range("I20").Value = 1.046386338933683
$ws.Range("J20").Value = 1.051241353772231
$ws.Range("K20").Value = 1.051668946592303
$ws.Range("L20").Value = 1.045939608883284
$ws.Range("M20").Value = 1.048470364215606
$ws.Range("N20").Value = 1.052734237073656
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.043545188967665
$ws.Range("D21").Value = 1.047066147106946
$ws.Range("E21").Value = 1.041210612750422
$ws.Range("F21").Value = 1.042253156902016
$ws.Range("I21").Value = 1.045713610627866
$ws.Range("J21").Value = 1.050077120011525
$ws.Range("K21").Value = 1.050605582463009
$ws.Range("L21").Value = 1.044771466370093
$ws.Range("M21").Value = 1.04581017765656
$ws.Range("N21").Value = 1.05156834996747
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.042560731591975
$ws.Range("D22").Value = 1.04626370341068
$ws.Range("E22").Value = 1.040340648134851
$ws.Range("F22").Value = 1.040435902780699
$ws.Range("I22").Value = 1.045287193368503
$ws.Range("J22").Value = 1.049341558403465
$ws.Range("K22").Value = 1.049933586697183
$ws.Range("L22").Value = 1.044033069003852
$ws.Range("M22").Value = 1.044127958948645
$ws.Range("N22").Value = 1.050831743777556
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.043082914227512
$ws.Range("D23").Value = 1.046689363898834
$ws.Range("E23").Value = 1.040802155586972
$ws.Range("F23").Value = 1.041400180101255
$ws.Range("I23").Value = 1.04551357466141
$ws.Range("J23").Value = 1.049731834747917
$ws.Range("K23").Value = 1.050290151593474
$ws.Range("L23").Value = 1.044424884545545
$ws.Range("M23").Value = 1.045020667726945
$ws.Range("N23").Value = 1.051222574359183
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.045132051877057
$ws.Range("D24").Value = 1.048359225811443
$ws.Range("E24").Value = 1.042611978099405
$ws.Range("F24").Value = 1.045176454873917
$ws.Range("I24").Value = 1.046397559505339
$ws.Range("J24").Value = 1.051260812731269
$ws.Range("K24").Value = 1.051686716865228
$ws.Range("L24").Value = 1.045959126976327
$ws.Range("M24").Value = 1.048514802868181
$ws.Range("N24").Value = 1.052753723666648
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.047496666397055
$ws.Range("D25").Value = 1.050285119375152
$ws.Range("E25").Value = 1.044697913954755
$ws.Range("F25").Value = 1.049519240356807
$ws.Range("I25").Value = 1.047408642672748
$ws.Range("J25").Value = 1.053019956657081
$ws.Range("K25").Value = 1.05329281121785
$ws.Range("L25").Value = 1.047722739703379
$ws.Range("M25").Value = 1.052529269098411
$ws.Range("N25").Value = 1.054515365778611

